$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.552.57"

$ws.Range("E2").Value = "  +0.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.509.49"

$ws.Range("E3").Value = "  -0.40%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.41"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -0.42%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.30"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +0.55%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.505.47"

$ws.Range("E7").Value = "  -0.49%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("D9").Style = "Normal"

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  +2.62%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.64"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  +8.37%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.45%  "

# Row 13 - Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.78"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  +2.63%  "

# Row 14 - ShibaInu
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000216"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  -1.83%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.105.21"

$ws.Range("E15").Value = "  -0.29%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.509.51"

$ws.Range("E16").Value = "  -0.52%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.752.48"

$ws.Range("E17").Value = "  +0.38%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  +0.17%  "

# Row 19 - Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.56"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +2.09%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +2.25%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.85"
$ws.Range("D21").Style = "Normal"

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "449.69"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  +0.60%  "

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.633"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  +1.47%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.651.16"

# Row 26 - row26 (PEPE -> Dai)
$ws.Range("B26").Value = "Dai"

$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  +0.00%  "

# Row 27 - row27 (Dai -> PEPE)
$ws.Range("B27").Value = "PEPE"

$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000127"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  -1.73%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.90"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  +6.34%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.09"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  -0.88%  "

# Row 30 - row30 (PancakeSwap -> Fetch.AI)
$ws.Range("B30").Value = "Fetch.AI"

$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.67"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  +6.59%  "

# Row 31 - row31 (Fetch.AI -> PancakeSwap)
$ws.Range("B31").Value = "PancakeSwap"

$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.52"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  +0.43%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  +1.62%  "

# Row 33 - Binance-PegBSC-USD
$ws.Range("E33").Value = "  +0.10%  "

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.78"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -0.12%  "

# Row 35 - NEARProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  +0.76%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +1.63%  "

# Row 37 - RenzoRestakedETH
$ws.Range("D37").Value = "3.504.65"

$ws.Range("E37").Value = "  -0.34%  "

# Row 38 - Aptos
$ws.Range("E38").Value = "  +0.08%  "

# Row 39 - USDe
$ws.Range("E39").Value = "  +0.04%  "

# Row 40 - Stacks
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.31"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  +5.16%  "

# Row 41 - FirstDigitalUSD
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  -0.11%  "

# Row 42 - row42 (Hedera -> Monero)
$ws.Range("B42").Value = "Monero"

$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "174.18"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  -1.81%  "

# Row 43 - row43 (Monero -> Hedera)
$ws.Range("B43").Value = "Hedera"

$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0899"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  +2.58%  "

# Row 44 - Filecoin
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.49"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  +0.86%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.37"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +11.77%  "

# Row 46 - Mantle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.882"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +0.15%  "

# Row 47 - OKB
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.92"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  +3.02%  "

# Row 48 - ONDO
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.30"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +3.80%  "

# Row 49 - row49 (dogwifhat -> Cosmos)
$ws.Range("B49").Value = "Cosmos"

$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.68"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  +1.20%  "

# Row 50 - row50 (Cosmos -> dogwifhat)
$ws.Range("B50").Value = "dogwifhat"

$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  -3.50%  "

# Row 51 - TheGraph
$ws.Range("E51").Value = "  +2.91%  "
